$d = $word.ActiveDocument

# --- First paragraph formatting -------------------------------------------------
# Add a thin paragraph border (5-twip gap on every side, default line style
# suppressed by only touching the distance) and widen the left indent from
# 120 to 225 twips (twips / 20 = points).
$p1 = $d.Paragraphs(1)
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromRight = 5
$p1.Range.ParagraphFormat.LeftIndent = 225 / 20

# --- First paragraph text ---------------------------------------------------
# The paragraph used to read "**ID__AFFARS_5319_topic_4__ID** " (placeholder
# run followed by a run containing a single trailing space). Locate the
# placeholder, drop the now-unwanted trailing-space run immediately after
# it, then rename the placeholder itself.
$oldId = "**ID__AFFARS_5319_topic_4__ID**"
$newId = "**ID__AFFARS_SUBPART_5319_5__ID**"

$found = $d.Content.Duplicate
$found.Find.Execute($oldId, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$trailingSpace = $d.Range($found.End, $found.End + 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

$d.Content.Find.Execute($oldId, $true, $false, $false, $false, $false, $true, 1, $false, $newId, 2)

Write-Output "done"
